# Update Name of Algo
# Apply corrected KNN-imputed values to the result data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.294
$ws.Range("A3").Value = -21.945
$ws.Range("A14").Value = -21.715
$ws.Range("A21").Value = -19.998
$ws.Range("A23").Value = -19.945
$ws.Range("A25").Value = -21.775
$ws.Range("C25").Value = -13.241
$ws.Range("A26").Value = -21.351
$ws.Range("C27").Value = -13.055
$ws.Range("A29").Value = -21.284
$ws.Range("C31").Value = -13.002
$ws.Range("C39").Value = -12.559
$ws.Range("C48").Value = -11.1
$ws.Range("C51").Value = -11.152
$ws.Range("C52").Value = -11.601
$ws.Range("A53").Value = -21.931
$ws.Range("C55").Value = -13.513
$ws.Range("C56").Value = -13.222
$ws.Range("A57").Value = -22.095
$ws.Range("C57").Value = -13.58
$ws.Range("A59").Value = -22.407
$ws.Range("A69").Value = -21.656
$ws.Range("C73").Value = -12.601
$ws.Range("A79").Value = -21.141
$ws.Range("A83").Value = -21.938
$ws.Range("C89").Value = -10.857
$ws.Range("C90").Value = -12.482
$ws.Range("A91").Value = -21.509
$ws.Range("C92").Value = -11.087
$ws.Range("A93").Value = -21.519

$wb.Save()
